# Generate Report for Handback
#
# Flips the per-language handoff status to "handed back" on the Overview
# sheet and fills in the newly-available handback info (target file,
# handback file + hyperlink, handback datetime) on the zh-cn / de-de
# per-language sheets.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both file rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $handedBack
$wsZhCn.Range("I2").Value = "03314ed2-a8fd-4811-9a3d-19d56888ab94.md"
$wsZhCn.Range("J2").Value = "03314ed2-a8fd-4811-9a3d-19d56888ab94.88f01aca8ed223ca81c52f73ab044c6c7ffe83bf.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-17 17:00:45"

$wsZhCn.Range("C3").Value = $handedBack
$wsZhCn.Range("I3").Value = "07f0619d-524a-4774-af47-6ae1dd83f36a.md"
$wsZhCn.Range("J3").Value = "07f0619d-524a-4774-af47-6ae1dd83f36a.2688b659d39c41a2c0c5f604915e2307dd9c1d8e.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-17 17:00:45"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c79329746a031f0e3a8c4daf7123687ed1283a57/e2e/03314ed2-a8fd-4811-9a3d-19d56888ab94.md", "", "", "03314ed2-a8fd-4811-9a3d-19d56888ab94.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c79329746a031f0e3a8c4daf7123687ed1283a57/e2e/07f0619d-524a-4774-af47-6ae1dd83f36a.md", "", "", "07f0619d-524a-4774-af47-6ae1dd83f36a.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $handedBack
$wsDeDe.Range("I2").Value = "03314ed2-a8fd-4811-9a3d-19d56888ab94.md"
$wsDeDe.Range("J2").Value = "03314ed2-a8fd-4811-9a3d-19d56888ab94.88f01aca8ed223ca81c52f73ab044c6c7ffe83bf.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-17 17:00:53"

$wsDeDe.Range("C3").Value = $handedBack
$wsDeDe.Range("I3").Value = "07f0619d-524a-4774-af47-6ae1dd83f36a.md"
$wsDeDe.Range("J3").Value = "07f0619d-524a-4774-af47-6ae1dd83f36a.2688b659d39c41a2c0c5f604915e2307dd9c1d8e.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-17 17:00:53"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c79329746a031f0e3a8c4daf7123687ed1283a57/e2e/03314ed2-a8fd-4811-9a3d-19d56888ab94.md", "", "", "03314ed2-a8fd-4811-9a3d-19d56888ab94.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c79329746a031f0e3a8c4daf7123687ed1283a57/e2e/07f0619d-524a-4774-af47-6ae1dd83f36a.md", "", "", "07f0619d-524a-4774-af47-6ae1dd83f36a.md")

# ---------------------------------------------------------------------
# Widen the (now longer) target/handback file columns so the new values
# are fully visible, matching the report generator's auto-fit behaviour.
# ---------------------------------------------------------------------
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

Write-Output "Generate Report for Handback: applied."
